$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("age")
$ws.Range("A548").Value = 44027
$ws.Range("C548").Value = 111
$ws.Range("D548").Value = 0.4225030451
$ws.Range("E548").Value = 5
$ws.Range("A549").Value = 44027
$ws.Range("C549").Value = 496
$ws.Range("D549").Value = 1.8879415347
$ws.Range("E549").Value = 19
$ws.Range("A550").Value = 44027
$ws.Range("C550").Value = 1293
$ws.Range("D550").Value = 4.921589525
$ws.Range("E550").Value = 49
$ws.Range("A551").Value = 44027
$ws.Range("C551").Value = 4901
$ws.Range("D551").Value = 18.654841657
$ws.Range("E551").Value = 191
$ws.Range("F551").Value = 4
$ws.Range("G551").Value = 0.5494505495
$ws.Range("H551").Value = -1
$ws.Range("A552").Value = 44027
$ws.Range("C552").Value = 4867
$ws.Range("D552").Value = 18.525426309
$ws.Range("E552").Value = 160
$ws.Range("G552").Value = 1.6483516484
$ws.Range("A553").Value = 44027
$ws.Range("C553").Value = 4684
$ws.Range("D553").Value = 17.828867235
$ws.Range("E553").Value = 152
$ws.Range("F553").Value = 28
$ws.Range("G553").Value = 3.8461538462
$ws.Range("H553").Value = 2
$ws.Range("A554").Value = 44027
$ws.Range("C554").Value = 4353
$ws.Range("D554").Value = 16.568970767
$ws.Range("E554").Value = 109
$ws.Range("G554").Value = 10.576923077
$ws.Range("A555").Value = 44027
$ws.Range("C555").Value = 1735
$ws.Range("D555").Value = 6.6039890378
$ws.Range("E555").Value = 47
$ws.Range("F555").Value = 74
$ws.Range("G555").Value = 10.164835165
$ws.Range("H555").Value = 2
$ws.Range("A556").Value = 44027
$ws.Range("C556").Value = 1266
$ws.Range("D556").Value = 4.818818514
$ws.Range("E556").Value = 40
$ws.Range("F556").Value = 84
$ws.Range("G556").Value = 11.538461538
$ws.Range("H556").Value = 2
$ws.Range("A557").Value = 44027
$ws.Range("C557").Value = 807
$ws.Range("D557").Value = 3.0717113276
$ws.Range("E557").Value = 11
$ws.Range("G557").Value = 10.851648352
$ws.Range("A558").Value = 44027
$ws.Range("C558").Value = 615
$ws.Range("D558").Value = 2.3408952497
$ws.Range("E558").Value = 14
$ws.Range("F558").Value = 85
$ws.Range("G558").Value = 11.675824176
$ws.Range("H558").Value = 2
$ws.Range("A559").Value = 44027
$ws.Range("C559").Value = 1126
$ws.Range("D559").Value = 4.2859317905
$ws.Range("E559").Value = 10
$ws.Range("F559").Value = 285
$ws.Range("G559").Value = 39.148351648
$ws.Range("H559").Value = 1
$ws.Range("A560").Value = 44027
$ws.Range("D560").Value = 0.0685140073

$ws = $wb.Worksheets.Item("gender")
$ws.Range("A128").Value = 44027
$ws.Range("C128").Value = 12560
$ws.Range("D128").Value = 47.807551766
$ws.Range("E128").Value = 333
$ws.Range("F128").Value = 300
$ws.Range("G128").Value = 41.208791209
$ws.Range("H128").Value = 3
$ws.Range("A129").Value = 44027
$ws.Range("C129").Value = 13369
$ws.Range("D129").Value = 50.886875761
$ws.Range("E129").Value = 424
$ws.Range("F129").Value = 427
$ws.Range("G129").Value = 58.653846154
$ws.Range("H129").Value = 7
$ws.Range("A130").Value = 44027
$ws.Range("C130").Value = 343
$ws.Range("D130").Value = 1.3055724726
$ws.Range("E130").Value = 50
$ws.Range("G130").Value = 0.1373626374

$ws = $wb.Worksheets.Item("race")
$ws.Range("A248").Value = 44027
$ws.Range("C248").Value = 651
$ws.Range("D248").Value = 2.4779232643
$ws.Range("E248").Value = 10
$ws.Range("G248").Value = 1.9230769231
$ws.Range("A249").Value = 44027
$ws.Range("C249").Value = 2989
$ws.Range("D249").Value = 11.377131547
$ws.Range("E249").Value = 74
$ws.Range("F249").Value = 94
$ws.Range("G249").Value = 12.912087912
$ws.Range("H249").Value = 2
$ws.Range("A250").Value = 44027
$ws.Range("C250").Value = 10534
$ws.Range("D250").Value = 40.09591961
$ws.Range("E250").Value = 424
$ws.Range("F250").Value = 211
$ws.Range("G250").Value = 28.983516484
$ws.Range("H250").Value = 1
$ws.Range("A251").Value = 44027
$ws.Range("C251").Value = 135
$ws.Range("D251").Value = 0.5138550548
$ws.Range("E251").Value = 2
$ws.Range("G251").Value = 0.1373626374
$ws.Range("A252").Value = 44027
$ws.Range("C252").Value = 5441
$ws.Range("D252").Value = 20.710261876
$ws.Range("E252").Value = 114
$ws.Range("F252").Value = 93
$ws.Range("G252").Value = 12.774725275
$ws.Range("H252").Value = 0
$ws.Range("A253").Value = 44027
$ws.Range("C253").Value = 6522
$ws.Range("D253").Value = 24.824908648
$ws.Range("E253").Value = 183
$ws.Range("F253").Value = 315
$ws.Range("G253").Value = 43.269230769
$ws.Range("H253").Value = 5
